# Maj Design Romain task
# Adds a second task block for Romain (rows 36-39), turns the old "box"
# separator rows (12, 19, 23, 24, 25, 26) into plain/blank separators
# (G:I, no border/fill), and gives H18 a 3-sided (no-bottom) red border
# so the box border continues visually into the new blank separator row 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# ---------------------------------------------------------------------
# 1) Strip the old full-box border from the "spacer" cells so they become
#    blank (no border / no fill) and extend the blank spacer across G:I.
# ---------------------------------------------------------------------
$spacerRanges = @("G12:I12", "G19:I19", "H23:H23", "G24:H24", "G25:I25", "G26:I26")
foreach ($rng in $spacerRanges) {
    $r = $ws.Range($rng)
    $r.Borders.LineStyle = -4142
}

# ---------------------------------------------------------------------
# 2) H18 keeps its red fill but now only has a border on the left/top/
#    right sides (no bottom), since the box continues into row 19.
# ---------------------------------------------------------------------
$h18 = $ws.Range("H18")
$h18.Borders.LineStyle = 1
$h18.Borders.Item(9).LineStyle = -4142

# ---------------------------------------------------------------------
# 3) New second task block for Romain (rows 36-39), mirroring the style
#    of the existing task blocks (yellow / red alternating box borders).
# ---------------------------------------------------------------------
$ws.Range("C36").Value = "Romain"

$ws.Range("G36").Value = "Design bumper (Generic)"
$ws.Range("H27").Copy($ws.Range("H36"))

$ws.Range("G37").Value = "Design platform (Generic)"
$ws.Range("H9").Copy($ws.Range("H37"))

$ws.Range("G38").Value = "Background music 1 (Paper)"
$ws.Range("H27").Copy($ws.Range("H38"))

$ws.Range("G39").Value = "Background music 2 (Paper)"
$ws.Range("H9").Copy($ws.Range("H39"))

# ---------------------------------------------------------------------
# 4) Scroll/selection bookkeeping to match the saved view state.
# ---------------------------------------------------------------------
$ws.Range("H34").Select()
